$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Tab-separated row data for the new "Vawa - Dating Violence" block (rows 3323-3394):
#   Row | A (Sector) | B (Reporting Location) | D (Date) | E (Count, blank = no cell) | HasEmptyStyledF (1/0)
$data = @"
3323	Public, 4-year or above	On Campus (excluding Residence Halls)	sum2014	616	0
3324	Private nonprofit, 4-year or above	On Campus (excluding Residence Halls)	sum2014	359	0
3325	Private for-profit, 4-year or above	On Campus (excluding Residence Halls)	sum2014	23	0
3326	Public, 2-year	On Campus (excluding Residence Halls)	sum2014	219	0
3327	Private nonprofit, 2-year	On Campus (excluding Residence Halls)	sum2014	4	0
3328	Private for-profit, 2-year	On Campus (excluding Residence Halls)	sum2014	12	0
3329	Public, less-than 2-year	On Campus (excluding Residence Halls)	sum2014	3	0
3330	Private nonprofit, less-than 2-year	On Campus (excluding Residence Halls)	sum2014	0	0
3331	Private for-profit, less-than 2-year	On Campus (excluding Residence Halls)	sum2014	1	0
3332	Public, 4-year or above	On Campus (excluding Residence Halls)	sum2015	703	0
3333	Private nonprofit, 4-year or above	On Campus (excluding Residence Halls)	sum2015	375	0
3334	Private for-profit, 4-year or above	On Campus (excluding Residence Halls)	sum2015	9	0
3335	Public, 2-year	On Campus (excluding Residence Halls)	sum2015	252	0
3336	Private nonprofit, 2-year	On Campus (excluding Residence Halls)	sum2015	1	0
3337	Private for-profit, 2-year	On Campus (excluding Residence Halls)	sum2015	9	0
3338	Public, less-than 2-year	On Campus (excluding Residence Halls)	sum2015	2	0
3339	Private nonprofit, less-than 2-year	On Campus (excluding Residence Halls)	sum2015	1	0
3340	Private for-profit, less-than 2-year	On Campus (excluding Residence Halls)	sum2015	6	0
3341	Public, 4-year or above	On Campus (Residence Halls)	sum2014	1147	1
3342	Private nonprofit, 4-year or above	On Campus (Residence Halls)	sum2014	737	0
3343	Private for-profit, 4-year or above	On Campus (Residence Halls)	sum2014	23	0
3344	Public, 2-year	On Campus (Residence Halls)	sum2014	88	0
3345	Private nonprofit, 2-year	On Campus (Residence Halls)	sum2014	5	0
3346	Private for-profit, 2-year	On Campus (Residence Halls)	sum2014	0	0
3347	Public, less-than 2-year	On Campus (Residence Halls)	sum2014	0	0
3348	Private nonprofit, less-than 2-year	On Campus (Residence Halls)	sum2014		0
3349	Private for-profit, less-than 2-year	On Campus (Residence Halls)	sum2014	0	0
3350	Public, 4-year or above	On Campus (Residence Halls)	sum2015	1521	0
3351	Private nonprofit, 4-year or above	On Campus (Residence Halls)	sum2015	958	0
3352	Private for-profit, 4-year or above	On Campus (Residence Halls)	sum2015	16	0
3353	Public, 2-year	On Campus (Residence Halls)	sum2015	144	0
3354	Private nonprofit, 2-year	On Campus (Residence Halls)	sum2015	3	0
3355	Private for-profit, 2-year	On Campus (Residence Halls)	sum2015	0	0
3356	Public, less-than 2-year	On Campus (Residence Halls)	sum2015	0	0
3357	Private nonprofit, less-than 2-year	On Campus (Residence Halls)	sum2015		0
3358	Private for-profit, less-than 2-year	On Campus (Residence Halls)	sum2015	1	0
3359	Public, 4-year or above	Non-Campus	sum2014	57	0
3360	Private nonprofit, 4-year or above	Non-Campus	sum2014	66	0
3361	Private for-profit, 4-year or above	Non-Campus	sum2014	3	0
3362	Public, 2-year	Non-Campus	sum2014	17	0
3363	Private nonprofit, 2-year	Non-Campus	sum2014	4	0
3364	Private for-profit, 2-year	Non-Campus	sum2014	0	0
3365	Public, less-than 2-year	Non-Campus	sum2014	0	0
3366	Private nonprofit, less-than 2-year	Non-Campus	sum2014	0	0
3367	Private for-profit, less-than 2-year	Non-Campus	sum2014	0	0
3368	Public, 4-year or above	Non-Campus	sum2015	72	0
3369	Private nonprofit, 4-year or above	Non-Campus	sum2015	42	0
3370	Private for-profit, 4-year or above	Non-Campus	sum2015	3	0
3371	Public, 2-year	Non-Campus	sum2015	23	0
3372	Private nonprofit, 2-year	Non-Campus	sum2015	0	0
3373	Private for-profit, 2-year	Non-Campus	sum2015	0	0
3374	Public, less-than 2-year	Non-Campus	sum2015	0	0
3375	Private nonprofit, less-than 2-year	Non-Campus	sum2015	0	0
3376	Private for-profit, less-than 2-year	Non-Campus	sum2015	1	0
3377	Public, 4-year or above	Public Property	sum2014	93	0
3378	Private nonprofit, 4-year or above	Public Property	sum2014	69	0
3379	Private for-profit, 4-year or above	Public Property	sum2014	4	0
3380	Public, 2-year	Public Property	sum2014	25	0
3381	Private nonprofit, 2-year	Public Property	sum2014	7	0
3382	Private for-profit, 2-year	Public Property	sum2014	8	0
3383	Public, less-than 2-year	Public Property	sum2014	2	0
3384	Private nonprofit, less-than 2-year	Public Property	sum2014	0	0
3385	Private for-profit, less-than 2-year	Public Property	sum2014	1	0
3386	Public, 4-year or above	Public Property	sum2015	81	0
3387	Private nonprofit, 4-year or above	Public Property	sum2015	74	0
3388	Private for-profit, 4-year or above	Public Property	sum2015	2	0
3389	Public, 2-year	Public Property	sum2015	28	0
3390	Private nonprofit, 2-year	Public Property	sum2015	0	0
3391	Private for-profit, 2-year	Public Property	sum2015	2	0
3392	Public, less-than 2-year	Public Property	sum2015	1	0
3393	Private nonprofit, less-than 2-year	Public Property	sum2015	1	0
3394	Private for-profit, less-than 2-year	Public Property	sum2015	5	0
"@

$offense = "Vawa - Dating Violence"

$lines = $data -split "`r?`n" | Where-Object { $_.Trim().Length -gt 0 }
foreach ($line in $lines) {
    $parts = $line -split "`t"
    $r      = [int]$parts[0]
    $aVal   = $parts[1]
    $bVal   = $parts[2]
    $dVal   = $parts[3]
    $eVal   = $parts[4]
    $hasF   = $parts[5]

    # Columns A, B, D carry the "quote-prefixed text" style (s="1") seen on the
    # rest of the sheet, so a leading apostrophe is used to force that style
    # while the apostrophe itself is not stored in the cell text.
    $ws.Cells.Item($r, 1).Value = "'" + $aVal
    $ws.Cells.Item($r, 2).Value = "'" + $bVal

    # Column C (Offense) is plain (unstyled) text.
    $ws.Cells.Item($r, 3).Value = $offense

    $ws.Cells.Item($r, 4).Value = "'" + $dVal

    # Column E (Count) is a plain number; some rows have no count at all.
    if ($eVal -ne "") {
        $ws.Cells.Item($r, 5).Value = [double]$eVal
    }

    # A couple of rows also carry a trailing empty-but-styled F cell, matching
    # the leftover formatting seen elsewhere in the sheet.
    if ($hasF -eq "1") {
        $ws.Cells.Item($r, 6).Value = "'"
        $ws.Cells.Item($r, 6).Value = ""
    }
}

# Update the view to match the edited state: scroll near the bottom of the
# sheet and leave the selection on the first empty row after the new data.
[void]$ws.Range("C3395").Select()
